$wb = $excel.ActiveWorkbook

# --- Sheet1: scroll back to top-left (drops the saved topLeftCell="A108"
#     view state) while preserving its existing selection A1:B1. This must
#     happen before we activate the new sheet below, since Range.Select()
#     also activates the sheet it is called on.
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Range("A1:B1").Select()

# --- Duplicate "20150812" into a new sheet "20150817" at the end of the
#     workbook (mirrors Excel's "Move or Copy... (Create a copy)").
$srcWs = $wb.Worksheets.Item("20150812")
$srcWs.Copy($null, $srcWs)
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "20150817"

# --- Cost update: refresh the values that changed between 2015-08-12 and
#     2015-08-17 (quintel/etsource#935).
$newWs.Range("B1").Value = 2723593.9226519298
$newWs.Range("B2").Value = 567345.39723509096
$newWs.Range("B3").Value = 1691150.70153317
$newWs.Range("B8").Value = 2723593.9226519298
$newWs.Range("B9").Value = 1691150.70153317
$newWs.Range("B10").Value = 567345.39723509096
$newWs.Range("B12").Value = 522910.85751158203
$newWs.Range("B13").Value = 2615893.9226519298
$newWs.Range("B14").Value = 1691150.70153317
$newWs.Range("B15").Value = 8542711.8949343897
$newWs.Range("B17").Value = 238299117.49039701
$newWs.Range("B20").Value = 123660642.97733299
$newWs.Range("B23").Value = 39843578.452615403
$newWs.Range("B27").Value = 279660189.01210499
$newWs.Range("B28").Value = 201505203.573333
$newWs.Range("B30").Value = 217173596.69085401
$newWs.Range("B31").Value = 112030560.895487
$newWs.Range("B32").Value = 76322643.130030796
$newWs.Range("B33").Value = 209387069.43603101
$newWs.Range("B35").Value = 234729173.82218
$newWs.Range("B36").Value = 215474700.58883801
$newWs.Range("B37").Value = 316926234.01324099
$newWs.Range("B38").Value = 182641746.32128
$newWs.Range("B39").Value = 26338032.819367401
$newWs.Range("B40").Value = 369101405.709297
$newWs.Range("B42").Value = 268611161.59955698
$newWs.Range("B44").Value = 1106589.20173967
$newWs.Range("B45").Value = 20655369.5223893
$newWs.Range("B46").Value = 44184731.0503067
$newWs.Range("B47").Value = 429600.48393909598
$newWs.Range("B48").Value = 5761794.4890970401
$newWs.Range("B49").Value = 190766827.816971
$newWs.Range("B54").Value = 585205.73772175598
$newWs.Range("B59").Value = 1971.0808485340899
$newWs.Range("B61").Value = 788303.766761926
$newWs.Range("B64").Value = 4460530.6897432096
$newWs.Range("B65").Value = 5831983.2552912803
$newWs.Range("B66").Value = 8802896.5365704894
$newWs.Range("B68").Value = 5831983.2552912803
$newWs.Range("B69").Value = 4121561.66890524
$newWs.Range("B70").Value = 5831983.2552912803
$newWs.Range("B71").Value = 4460530.6897432096
$newWs.Range("B73").Value = 8802896.5365704894
$newWs.Range("B74").Value = 5831983.2552912803
$newWs.Range("B75").Value = 3923656.88004041
$newWs.Range("B80").Value = 155.96499963094899
$newWs.Range("B83").Value = 290.546981206743
$newWs.Range("B91").Value = 444.98040974133102
$newWs.Range("B94").Value = 413.427361255606
$newWs.Range("B99").Value = 4460530.6897432096
$newWs.Range("B102").Value = 288730.04206686502
$newWs.Range("B103").Value = 288730.04206686502
$newWs.Range("B108").Value = 4684910.64292455
$newWs.Range("B109").Value = 6919979.6348137204
$newWs.Range("B111").Value = 5831983.2552912803
$newWs.Range("B115").Value = 10025561.2852132

# --- Formatting tweaks: drop the best-fit marker on column B and add a
#     blank, standard-width column C; set the new sheet's print setup.
$newWs.Columns.Item(2).ColumnWidth = $newWs.Columns.Item(2).ColumnWidth
$newWs.Columns.Item(3).ColumnWidth = 10
$newWs.PageSetup.PaperSize = 9
$newWs.PageSetup.Orientation = 1

# --- Old "20150812" sheet is no longer the tab in focus: narrow its saved
#     selection to column A and give it a print setup too.
[void]$srcWs.Range("A1:A121").Select()
$srcWs.PageSetup.PaperSize = 9
$srcWs.PageSetup.Orientation = 1

# --- The newly added sheet becomes the active / focused tab.
[void]$newWs.Range("A1").Select()
$newWs.Activate()
